$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = -7.402699999999993
$ws.Range("C7").Value = -13.12159999999999
$ws.Range("A8").Value = -22.30200000000001
$ws.Range("A10").Value = -21.79579999999999
$ws.Range("A12").Value = -21.5587
$ws.Range("C15").Value = -14.37019999999998
$ws.Range("A18").Value = -22.0209
$ws.Range("C18").Value = -12.87009999999999
$ws.Range("D18").Value = -8.280499999999991
$ws.Range("D19").Value = -8.923399999999996
$ws.Range("C20").Value = -12.20720000000001
$ws.Range("D27").Value = -8.624400000000005
$ws.Range("C29").Value = -11.5497
$ws.Range("C30").Value = -12.8468
$ws.Range("C31").Value = -12.7335
$ws.Range("D31").Value = -8.950300000000004
$ws.Range("A37").Value = -20.88660000000001
$ws.Range("D38").Value = -8.443199999999997
$ws.Range("C40").Value = -13.3159
$ws.Range("D42").Value = -8.907699999999995
$ws.Range("D44").Value = -7.504099999999998
$ws.Range("D47").Value = -7.620600000000001
$ws.Range("C50").Value = -13.69799999999999
$ws.Range("A55").Value = -22.37030000000001
$ws.Range("D58").Value = -8.486899999999993
$ws.Range("D65").Value = -7.658099999999999
$ws.Range("A68").Value = -21.70579999999999
$ws.Range("C68").Value = -12.3008
$ws.Range("D73").Value = -7.843199999999997
$ws.Range("C76").Value = -12.44359999999999
$ws.Range("A77").Value = -20.9585
$ws.Range("A78").Value = -20.53809999999999
$ws.Range("A81").Value = -21.86250000000001
$ws.Range("A82").Value = -21.8471
$ws.Range("C87").Value = -13.53749999999999
$ws.Range("C88").Value = -13.22999999999999
$ws.Range("D90").Value = -7.970700000000004
$ws.Range("D94").Value = -6.9429
$ws.Range("D95").Value = -7.707199999999999
$ws.Range("C96").Value = -13.15500000000001
$ws.Range("C98").Value = -12.24189999999999
$ws.Range("C101").Value = -13.6473
$ws.Range("D101").Value = -8.176499999999992
$ws.Range("C102").Value = -13.29330000000001
